$wb = $excel.ActiveWorkbook

# --- Sheet 1: ASSESSMENTS (text stays, only the instructions text changes) ---
$ws1 = $wb.Worksheets.Item("ASSESSMENTS")
$ws1.Range("E5").Value = "Add the assessments you want to add to the tool on the first sheet.`n If the name of any assessment coincides with an existing one, the latter (along with its sites) will be deleted."
$ws1.Range("E6").Value = "In the sites tab, define the industries you want to add,  `nwith the assessment to which it belongs. This assessment does not have to be defined in the assessment sheet, it can be previously defined in the web tool."

# --- Sheet 2: INDUSTRIES -> SITES ---
$ws2 = $wb.Worksheets.Item("INDUSTRIES")
$ws2.Name = "SITES"

$ws2.Range("A1").Value = "SITE"
$ws2.Range("E1").Value = "SUB-SUPPLIERS"
$ws2.Range("E2").Value = "Add as many sub-suppliers (Name, Latitude, Longitude) as needed to the right"
$ws2.Range("B3").Value = "Site "

$ws1.Range("E13").Select()
$ws1.Range("E6").Select()

$ws2.Range("J13").Select()
$ws2.Range("B28").Select()
